# Swap the data values (excluding "Rödlistade" column D and "Lokalnamn" column P,
# which are identical between the paired rows) between row 3 <-> row 4, and
# between row 20 <-> row 21, on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("A", "B", "E", "F", "G", "H", "Q", "R", "AC")

function Swap-RowValues($ws, $row1, $row2, $columns) {
    foreach ($col in $columns) {
        $addr1 = "$col$row1"
        $addr2 = "$col$row2"
        $val1 = $ws.Range($addr1).Value2
        $val2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value2 = $val2
        $ws.Range($addr2).Value2 = $val1
    }
}

Swap-RowValues $ws 3 4 $columns
Swap-RowValues $ws 20 21 $columns
